# SwaadSutra_Daily_2026-01-20.xlsx
# New order #22 ("Til Poli x1", Pooja, flat 12) came in and was inserted at the
# top of the "Daily Orders" log (newest-first), pushing the existing rows down.
# The "Summary" and "Items Breakdown" sheets are recomputed to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Daily Orders" - insert the new order as row 2 (pushes old rows down)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daily Orders")
$ws1.Rows.Item(2).Insert()

$ws1.Cells.Item(2,1).Value  = 22
$ws1.Cells.Item(2,2).Value  = "2026-01-20 11:13"
$ws1.Cells.Item(2,3).Value  = "Pooja"
# Flat No / Phone / Collection Date look numeric - force them to stay text
# (quote-prefix), matching the rest of the column's text-typed cells.
$ws1.Cells.Item(2,4).Value  = "'12"
$ws1.Cells.Item(2,5).Value  = "'9096648553"
$ws1.Cells.Item(2,6).Value  = "Til Poli x1"
$ws1.Cells.Item(2,7).Value  = 30
$ws1.Cells.Item(2,8).Value  = "NEW"
$ws1.Cells.Item(2,9).Value  = "PENDING"
$ws1.Cells.Item(2,10).Value = "'2026-01-20"
$ws1.Cells.Item(2,11).Value = "16:43"
$ws1.Cells.Item(2,12).Value = "'"
$ws1.Cells.Item(2,13).Value = "'"
$ws1.Cells.Item(2,14).Value = "'"

# ---------------------------------------------------------------------------
# Sheet 2: "Summary" - bump Total Orders, New count, and Total Revenue
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Cells.Item(2,1).Value = 3   # Total Orders: 2 -> 3
$ws2.Cells.Item(2,2).Value = 1   # New: 0 -> 1
$ws2.Cells.Item(2,7).Value = 80  # Total Revenue: 50 -> 80

# ---------------------------------------------------------------------------
# Sheet 3: "Items Breakdown" - add "Til Poli" line (insert above Jawar Bhakari)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Items Breakdown")
$ws3.Rows.Item(3).Insert()
$ws3.Cells.Item(3,1).Value = "Til Poli"
$ws3.Cells.Item(3,2).Value = 1
$ws3.Cells.Item(3,3).Value = 30
